# Scheduled runner update: refresh currentAveragePrice / Leve profit
# columns across all class sheets with the latest Universalis snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 90.333336
$ws.Range("I6").Value = 88.40000000000001
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 265.2
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -153.2
$ws.Range("N6").Value = -524

$ws.Range("H111").Value = 4649.1665
$ws.Range("J111").Value = 4874
$ws.Range("L111").Value = 14622
$ws.Range("N111").Value = -20756

$ws.Range("H113").Value = 2999.5
$ws.Range("I113").Value = 2999.5
$ws.Range("K113").Value = 2999.5
$ws.Range("M113").Value = 254.5

$ws.Range("H116").Value = 12456.615
$ws.Range("I116").Value = 14290.059
$ws.Range("J116").Value = 8993.444
$ws.Range("K116").Value = 14290.059
$ws.Range("L116").Value = 8993.444
$ws.Range("M116").Value = -10848.059
$ws.Range("N116").Value = -15877.444

$ws.Range("H118").Value = 2031.1428
$ws.Range("I118").Value = 2043.8
$ws.Range("K118").Value = 6131.4
$ws.Range("M118").Value = -4474.4

$ws.Range("H132").Value = 86093.53999999999
$ws.Range("I132").Value = 95732.56
$ws.Range("K132").Value = 287197.68
$ws.Range("M132").Value = -284667.68

$ws.Range("H138").Value = 2809
$ws.Range("I138").Value = 2271.5557
$ws.Range("J138").Value = 4421.3335
$ws.Range("K138").Value = 6814.6671
$ws.Range("L138").Value = 13264.0005
$ws.Range("M138").Value = -1674.6671
$ws.Range("N138").Value = -23544.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1827.766
$ws.Range("I2").Value = 1644.4318
$ws.Range("K2").Value = 1644.4318
$ws.Range("M2").Value = -1531.4318

$ws.Range("H31").Value = 26515
$ws.Range("I31").Value = 26515
$ws.Range("K31").Value = 26515
$ws.Range("M31").Value = -26221

$ws.Range("H32").Value = 4731.18
$ws.Range("I32").Value = 2545.3252
$ws.Range("J32").Value = 15403.294
$ws.Range("K32").Value = 2545.3252
$ws.Range("L32").Value = 15403.294
$ws.Range("M32").Value = -2258.3252
$ws.Range("N32").Value = -15977.294

$ws.Range("H61").Value = 3047895.8
$ws.Range("I61").Value = 3047895.8
$ws.Range("K61").Value = 3047895.8
$ws.Range("M61").Value = -3047683.8

$ws.Range("H92").Value = 77685.71000000001
$ws.Range("J92").Value = 77685.71000000001
$ws.Range("L92").Value = 77685.71000000001
$ws.Range("N92").Value = -82677.71000000001

$ws.Range("H97").Value = 892.1905
$ws.Range("I97").Value = 958.2222
$ws.Range("J97").Value = 496
$ws.Range("K97").Value = 958.2222
$ws.Range("L97").Value = 496
$ws.Range("M97").Value = -462.2222
$ws.Range("N97").Value = -1488

$ws.Range("H102").Value = 3967.9412
$ws.Range("I102").Value = 2904.6428
$ws.Range("K102").Value = 2904.6428
$ws.Range("M102").Value = -1282.6428

$ws.Range("H116").Value = 1827.766
$ws.Range("I116").Value = 1644.4318
$ws.Range("K116").Value = 1644.4318
$ws.Range("M116").Value = 649.5681999999999

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 2682.48
$ws.Range("I122").Value = 1156.0769
$ws.Range("J122").Value = 4336.0835
$ws.Range("K122").Value = 3468.2307
$ws.Range("L122").Value = 13008.2505
$ws.Range("M122").Value = -1018.2307
$ws.Range("N122").Value = -17908.2505

$ws.Range("H132").Value = 1228883.5
$ws.Range("I132").Value = 1553466
$ws.Range("K132").Value = 4660398
$ws.Range("M132").Value = -4657868

$ws.Range("H136").Value = 3047895.8
$ws.Range("I136").Value = 3047895.8
$ws.Range("K136").Value = 9143687.399999999
$ws.Range("M136").Value = -9141137.399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1827.766
$ws.Range("I3").Value = 1644.4318
$ws.Range("K3").Value = 1644.4318
$ws.Range("M3").Value = -1530.4318

$ws.Range("H99").Value = 15565.25
$ws.Range("I99").Value = 17878.29
$ws.Range("K99").Value = 17878.29
$ws.Range("M99").Value = -16380.29

$ws.Range("H105").Value = 1784.15
$ws.Range("I105").Value = 1823.8667
$ws.Range("K105").Value = 1823.8667
$ws.Range("M105").Value = -76.86670000000004

$ws.Range("H107").Value = 7863
$ws.Range("I107").Value = 7863
$ws.Range("K107").Value = 7863
$ws.Range("M107").Value = -5943

$ws.Range("H135").Value = 99984.5
$ws.Range("J135").Value = 99984.5
$ws.Range("L135").Value = 99984.5
$ws.Range("N135").Value = -110124.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 917
$ws.Range("I16").Value = 842.8
$ws.Range("K16").Value = 842.8
$ws.Range("M16").Value = -555.8

$ws.Range("H99").Value = 4202.161
$ws.Range("I99").Value = 2946.6667
$ws.Range("K99").Value = 2946.6667
$ws.Range("M99").Value = -1448.6667

$ws.Range("H107").Value = 958.8182
$ws.Range("I107").Value = 1021.8889
$ws.Range("J107").Value = 675
$ws.Range("K107").Value = 1021.8889
$ws.Range("L107").Value = 675
$ws.Range("M107").Value = 898.1111
$ws.Range("N107").Value = -4515

$ws.Range("H113").Value = 917
$ws.Range("I113").Value = 842.8
$ws.Range("K113").Value = 842.8
$ws.Range("M113").Value = 1327.2

$ws.Range("H126").Value = 4202.161
$ws.Range("I126").Value = 2946.6667
$ws.Range("K126").Value = 8840.000100000001
$ws.Range("M126").Value = -6370.000100000001

$ws.Range("H132").Value = 10016263
$ws.Range("I132").Value = 18979.62
$ws.Range("K132").Value = 56938.86
$ws.Range("M132").Value = -54408.86

$ws.Range("H141").Value = 216966.94
$ws.Range("J141").Value = 216966.94
$ws.Range("L141").Value = 216966.94
$ws.Range("N141").Value = -227326.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 3154087.8
$ws.Range("I7").Value = 2500323.8
$ws.Range("K7").Value = 7500971.399999999
$ws.Range("M7").Value = -7500859.399999999

$ws.Range("H92").Value = 981.8946999999999
$ws.Range("J92").Value = 1219.75
$ws.Range("L92").Value = 3659.25
$ws.Range("N92").Value = -6155.25

$ws.Range("H107").Value = 907.8
$ws.Range("I107").Value = 629
$ws.Range("J107").Value = 1093.6666
$ws.Range("K107").Value = 1887
$ws.Range("L107").Value = 3280.9998
$ws.Range("M107").Value = 33
$ws.Range("N107").Value = -7120.9998

$ws.Range("H113").Value = 1503.9048
$ws.Range("J113").Value = 1343.3529
$ws.Range("L113").Value = 4030.0587
$ws.Range("N113").Value = -8370.058700000001

$ws.Range("H140").Value = 3759.6667
$ws.Range("I140").Value = 3192.3635
$ws.Range("K140").Value = 9577.0905
$ws.Range("M140").Value = -4397.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7697.9
$ws.Range("I80").Value = 7087.636
$ws.Range("K80").Value = 7087.636
$ws.Range("M80").Value = -6089.636

$ws.Range("H83").Value = 7697.9
$ws.Range("I83").Value = 7087.636
$ws.Range("K83").Value = 35438.18
$ws.Range("M83").Value = -30446.18

$ws.Range("H135").Value = 45793.6
$ws.Range("J135").Value = 45793.6
$ws.Range("L135").Value = 45793.6
$ws.Range("N135").Value = -55933.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8894.583000000001
$ws.Range("I16").Value = 8894.583000000001
$ws.Range("K16").Value = 8894.583000000001
$ws.Range("M16").Value = -8724.583000000001

$ws.Range("H93").Value = 2380.8333
$ws.Range("I93").Value = 2596.5
$ws.Range("J93").Value = 2273
$ws.Range("K93").Value = 2596.5
$ws.Range("L93").Value = 2273
$ws.Range("M93").Value = -1348.5
$ws.Range("N93").Value = -4769

$ws.Range("H132").Value = 754609.5600000001
$ws.Range("J132").Value = 5295.1816
$ws.Range("L132").Value = 15885.5448
$ws.Range("N132").Value = -20945.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2998.158
$ws.Range("I107").Value = 1292.4
$ws.Range("K107").Value = 3877.2
$ws.Range("M107").Value = -1957.2

$ws.Range("H109").Value = 79000
$ws.Range("J109").Value = 79000
$ws.Range("L109").Value = 79000
$ws.Range("N109").Value = -81774
